# Apply the CodeSystem-gender-identity metadata update:
#  - Version bumped 5.0.0 -> 6.0.0
#  - Date bumped 2021-12-16T17:36:56+00:00 -> 2022-01-21T20:46:54+00:00
#  - Publisher now has a value: "Alvearie Team"
#  - The old "Contact" / "No display for ContactDetail" rows are replaced by a
#    "Jurisdiction" / "United States of America" row
#  - Every metadata row from "Description" through "Count" shifts up one row
#    (the old row 22 "Count" is absorbed and the very last row is removed)
#  - "Case Sensitive" now has the value "true"
#  - The final (old) row 22 is removed entirely, shrinking the sheet to A1:B21

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# --- Update values that stay in place ---
$ws.Range("B3").Value = "6.0.0"
$ws.Range("B8").Value = "2022-01-21T20:46:54+00:00"
$ws.Range("B9").Value = "Alvearie Team"

# --- Row 10: was Contact / No display for ContactDetail -> Jurisdiction / United States of America ---
$ws.Range("A10").Value = "Jurisdiction"
$ws.Range("B10").Value = "United States of America"

# --- Row 11: was Contact / No display for ContactDetail -> Description / <long description text> ---
$ws.Range("A11").Value = "Description"
$ws.Range("B11").Value = "Typically, gender identity codes are described in terms of masculinity and femininity. Gender identity is a social construct that varies, both across different cultures and over time. This value is often used for identity purposes, and should be collected directly from the patient."

# --- Row 12: was Description / <long text> -> Purpose / (empty) ---
$ws.Range("A12").Value = "Purpose"
$ws.Range("B12").Value = ""

# --- Row 13: was Purpose / (empty) -> Copyright / (empty) ---
$ws.Range("A13").Value = "Copyright"
$ws.Range("B13").Value = ""

# --- Row 14: was Copyright / (empty) -> Case Sensitive / true ---
$ws.Range("A14").Value = "Case Sensitive"
# Force the word "true" to be stored as literal text, not a boolean:
# write it with a leading apostrophe (Excel's "force text" marker), then
# restore the cell's normal (non quote-prefixed) formatting by pasting just
# the number format from a neighboring, untouched cell that already carries
# the correct style.
$ws.Range("B14").Value = "'true"
$ws.Range("B15").Copy() | Out-Null
$ws.Range("B14").PasteSpecial(-4122) | Out-Null

# --- Row 15: was Case Sensitive / (empty) -> Value Set (all codes) / (empty) ---
$ws.Range("A15").Value = "Value Set (all codes)"
$ws.Range("B15").Value = ""

# --- Row 16: was Value Set (all codes) / (empty) -> Hierarchy / (empty) ---
$ws.Range("A16").Value = "Hierarchy"
$ws.Range("B16").Value = ""

# --- Row 17: was Hierarchy / (empty) -> Compositional / (empty) ---
$ws.Range("A17").Value = "Compositional"
$ws.Range("B17").Value = ""

# --- Row 18: was Compositional / (empty) -> Version Needed? / (empty) ---
$ws.Range("A18").Value = "Version Needed?"
$ws.Range("B18").Value = ""

# --- Row 19: was Version Needed? / (empty) -> Content / complete ---
$ws.Range("A19").Value = "Content"
$ws.Range("B19").Value = "complete"

# --- Row 20: was Content / complete -> Supplements / (empty) ---
$ws.Range("A20").Value = "Supplements"
$ws.Range("B20").Value = ""

# --- Row 21: was Supplements / (empty) -> Count / 6 ---
$ws.Range("A21").Value = "Count"
# Force "6" to be stored as literal text, not a number, using the same
# apostrophe + format-paste trick as above.
$ws.Range("B21").Value = "'6"
$ws.Range("B20").Copy() | Out-Null
$ws.Range("B21").PasteSpecial(-4122) | Out-Null

# --- Row 22 (old "Count" / "6") is no longer needed; remove it entirely ---
$ws.Rows.Item(22).Delete()
